$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")

# ALC row 33
$ws_ALC.Range("H33").Value = 28606406
$ws_ALC.Range("I33").Value = 40000856
$ws_ALC.Range("J33").Value = 120277
$ws_ALC.Range("K33").Value = 40000856
$ws_ALC.Range("L33").Value = 120277
$ws_ALC.Range("M33").Value = -40000627
$ws_ALC.Range("N33").Value = -120735

# ALC row 40
$ws_ALC.Range("H40").Value = 10536.154
$ws_ALC.Range("I40").Value = 26142.5
$ws_ALC.Range("J40").Value = 3600
$ws_ALC.Range("K40").Value = 26142.5
$ws_ALC.Range("L40").Value = 3600
$ws_ALC.Range("M40").Value = -25967.5
$ws_ALC.Range("N40").Value = -3950

# ALC row 64
$ws_ALC.Range("H64").Value = 3358.652
$ws_ALC.Range("I64").Value = 3297.7778
$ws_ALC.Range("J64").Value = 3397.7856
$ws_ALC.Range("K64").Value = 3297.7778
$ws_ALC.Range("L64").Value = 3397.7856
$ws_ALC.Range("M64").Value = -3049.7778
$ws_ALC.Range("N64").Value = -3893.7856

# ALC row 67
$ws_ALC.Range("H67").Value = 3358.652
$ws_ALC.Range("I67").Value = 3297.7778
$ws_ALC.Range("J67").Value = 3397.7856
$ws_ALC.Range("K67").Value = 3297.7778
$ws_ALC.Range("L67").Value = 3397.7856
$ws_ALC.Range("M67").Value = -2439.7778
$ws_ALC.Range("N67").Value = -5113.7856

# ALC row 129
$ws_ALC.Range("H129").Value = 3220.7144

# ALC row 137
$ws_ALC.Range("H137").Value = 972.4545000000001
$ws_ALC.Range("J137").Value = 1233.6666
$ws_ALC.Range("L137").Value = 3700.9998
$ws_ALC.Range("N137").Value = -8800.9998

# ARM row 61
$ws_ARM.Range("H61").Value = 941.3509
$ws_ARM.Range("J61").Value = 1200
$ws_ARM.Range("L61").Value = 1200
$ws_ARM.Range("N61").Value = -1624

# ARM row 74
$ws_ARM.Range("H74").Value = 820.3889
$ws_ARM.Range("I74").Value = 670
$ws_ARM.Range("J74").Value = 1121.1666
$ws_ARM.Range("K74").Value = 670
$ws_ARM.Range("L74").Value = 1121.1666
$ws_ARM.Range("M74").Value = 204
$ws_ARM.Range("N74").Value = -2869.1666

# ARM row 77
$ws_ARM.Range("H77").Value = 820.3889
$ws_ARM.Range("I77").Value = 670
$ws_ARM.Range("J77").Value = 1121.1666
$ws_ARM.Range("K77").Value = 3350
$ws_ARM.Range("L77").Value = 5605.833000000001
$ws_ARM.Range("M77").Value = 1018
$ws_ARM.Range("N77").Value = -14341.833

# ARM row 88
$ws_ARM.Range("H88").Value = 1783.3334
$ws_ARM.Range("I88").Value = 1450
$ws_ARM.Range("J88").Value = 1950
$ws_ARM.Range("K88").Value = 1450
$ws_ARM.Range("L88").Value = 1950
$ws_ARM.Range("M88").Value = -1044
$ws_ARM.Range("N88").Value = -2762

# ARM row 91
$ws_ARM.Range("H91").Value = 1783.3334
$ws_ARM.Range("I91").Value = 1450
$ws_ARM.Range("J91").Value = 1950
$ws_ARM.Range("K91").Value = 1450
$ws_ARM.Range("L91").Value = 1950
$ws_ARM.Range("M91").Value = -46
$ws_ARM.Range("N91").Value = -4758

# ARM row 110
$ws_ARM.Range("H110").Value = 1647.6154
$ws_ARM.Range("I110").Value = 891.9
$ws_ARM.Range("J110").Value = 4166.6665
$ws_ARM.Range("K110").Value = 891.9
$ws_ARM.Range("L110").Value = 4166.6665
$ws_ARM.Range("M110").Value = 1153.1
$ws_ARM.Range("N110").Value = -8256.666499999999

# ARM row 125
$ws_ARM.Range("H125").Value = 17999.5
$ws_ARM.Range("J125").Value = 17999.5
$ws_ARM.Range("L125").Value = 17999.5
$ws_ARM.Range("N125").Value = -27839.5

# ARM row 132
$ws_ARM.Range("H132").Value = 1902.2
$ws_ARM.Range("I132").Value = 1545.75
$ws_ARM.Range("J132").Value = 3328
$ws_ARM.Range("K132").Value = 4637.25
$ws_ARM.Range("L132").Value = 9984
$ws_ARM.Range("M132").Value = -2107.25
$ws_ARM.Range("N132").Value = -15044

# ARM row 136
$ws_ARM.Range("H136").Value = 941.3509
$ws_ARM.Range("J136").Value = 1200
$ws_ARM.Range("L136").Value = 3600
$ws_ARM.Range("N136").Value = -8700

# BSM row 86
$ws_BSM.Range("H86").Value = 2128.5715
$ws_BSM.Range("I86").Value = 1500
$ws_BSM.Range("J86").Value = 2380
$ws_BSM.Range("K86").Value = 1500
$ws_BSM.Range("L86").Value = 2380
$ws_BSM.Range("M86").Value = -377
$ws_BSM.Range("N86").Value = -4626

# BSM row 89
$ws_BSM.Range("H89").Value = 2128.5715
$ws_BSM.Range("I89").Value = 1500
$ws_BSM.Range("J89").Value = 2380
$ws_BSM.Range("K89").Value = 7500
$ws_BSM.Range("L89").Value = 11900
$ws_BSM.Range("M89").Value = -1884
$ws_BSM.Range("N89").Value = -23132

# BSM row 94
$ws_BSM.Range("H94").Value = 371.72726
$ws_BSM.Range("I94").Value = 349
$ws_BSM.Range("J94").Value = 599
$ws_BSM.Range("K94").Value = 349
$ws_BSM.Range("L94").Value = 599
$ws_BSM.Range("M94").Value = 102
$ws_BSM.Range("N94").Value = -1501

# BSM row 105
$ws_BSM.Range("H105").Value = 2068158.4
$ws_BSM.Range("I105").Value = 2842717.8
$ws_BSM.Range("J105").Value = 2666.6667
$ws_BSM.Range("K105").Value = 2842717.8
$ws_BSM.Range("L105").Value = 2666.6667
$ws_BSM.Range("M105").Value = -2840970.8
$ws_BSM.Range("N105").Value = -6160.6667

# CRP row 7
$ws_CRP.Range("H7").Value = 124.1
$ws_CRP.Range("I7").Value = 130.125
$ws_CRP.Range("K7").Value = 130.125
$ws_CRP.Range("M7").Value = -17.125

# CRP row 62
$ws_CRP.Range("H62").Value = 3739.7
$ws_CRP.Range("I62").Value = 3628.1428
$ws_CRP.Range("K62").Value = 3628.1428
$ws_CRP.Range("M62").Value = -3004.1428

# CRP row 65
$ws_CRP.Range("H65").Value = 3739.7
$ws_CRP.Range("I65").Value = 3628.1428
$ws_CRP.Range("K65").Value = 18140.714
$ws_CRP.Range("M65").Value = -15020.714

# CUL row 70
$ws_CUL.Range("H70").Value = 3844.5715
$ws_CUL.Range("I70").Value = 956
$ws_CUL.Range("J70").Value = 5000
$ws_CUL.Range("K70").Value = 2868
$ws_CUL.Range("L70").Value = 15000
$ws_CUL.Range("M70").Value = -2553
$ws_CUL.Range("N70").Value = -15630

# CUL row 73
$ws_CUL.Range("H73").Value = 3844.5715
$ws_CUL.Range("I73").Value = 956
$ws_CUL.Range("J73").Value = 5000
$ws_CUL.Range("K73").Value = 2868
$ws_CUL.Range("L73").Value = 15000
$ws_CUL.Range("M73").Value = -1776
$ws_CUL.Range("N73").Value = -17184

# CUL row 133
$ws_CUL.Range("H133").Value = 7686.25
$ws_CUL.Range("J133").Value = 8583.333000000001
$ws_CUL.Range("L133").Value = 25749.999
$ws_CUL.Range("N133").Value = -35869.999

# GSM row 80
$ws_GSM.Range("H80").Value = 2133.889
$ws_GSM.Range("I80").Value = 2002.5
$ws_GSM.Range("J80").Value = 2171.4285
$ws_GSM.Range("K80").Value = 2002.5
$ws_GSM.Range("L80").Value = 2171.4285
$ws_GSM.Range("M80").Value = -1004.5
$ws_GSM.Range("N80").Value = -4167.4285

# GSM row 83
$ws_GSM.Range("H83").Value = 2133.889
$ws_GSM.Range("I83").Value = 2002.5
$ws_GSM.Range("J83").Value = 2171.4285
$ws_GSM.Range("K83").Value = 10012.5
$ws_GSM.Range("L83").Value = 10857.1425
$ws_GSM.Range("M83").Value = -5020.5
$ws_GSM.Range("N83").Value = -20841.1425

# LTW row 68
$ws_LTW.Range("H68").Value = 5643.75
$ws_LTW.Range("I68").Value = 8213.571
$ws_LTW.Range("J68").Value = 2046
$ws_LTW.Range("K68").Value = 8213.571
$ws_LTW.Range("L68").Value = 2046
$ws_LTW.Range("M68").Value = -7464.571
$ws_LTW.Range("N68").Value = -3544

# LTW row 71
$ws_LTW.Range("H71").Value = 5643.75
$ws_LTW.Range("I71").Value = 8213.571
$ws_LTW.Range("J71").Value = 2046
$ws_LTW.Range("K71").Value = 41067.855
$ws_LTW.Range("L71").Value = 10230
$ws_LTW.Range("M71").Value = -37323.855
$ws_LTW.Range("N71").Value = -17718

